$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.891.04'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.639.70'
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.88'
$ws.Range("E5").Value = '  +2.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '323.10'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.541'
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.68'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.77'
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.22'
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").Value = '3.040.49'
$ws.Range("E15").Value = '  +3.51%  '
$ws.Range("D16").Value = '2.632.08'
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.865'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '48.863.46'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.85'
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.68'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("D22").Value = '0.0₃0941'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.96'
$ws.Range("E23").Value = '  -4.29%  '
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.10'
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("E28").Value = '  +3.71%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.99'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").Value = '  -4.85%  '
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.43'
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.28'
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0798'
$ws.Range("E36").Value = '  +2.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.92'
$ws.Range("E37").Value = '  +5.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("E39").Value = '  +6.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.13'
$ws.Range("E40").Value = '  +4.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.53'
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("E43").Value = '  -3.49%  '
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").Value = '2.068.53'
$ws.Range("E45").Value = '  +2.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.22'
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  +6.36%  '
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.95'
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("E51").Value = '  -1.90%  '
